$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "survey" (sheet 1): insert a new "appearance" column after
# "type", and add a second "yes_no" select_one group (grid/inline
# appearance properties).
# ---------------------------------------------------------------
$survey = $wb.Worksheets.Item(1)

# Insert a new column C ("appearance"), shifting condition/name/label
# (old C,D,E) to D,E,F. Column widths tag along automatically.
$survey.Columns.Item(3).Insert()

$survey.Range("C1").Value = "appearance"
$survey.Range("C5").Value = "grid"

$survey.Range("B9").Value = "begin screen"

$survey.Range("B10").Value = "select_one yes_no"
$survey.Range("C10").Value = "inline"
$survey.Range("E10").Value = "i1"
$survey.Range("F10").Value = "Choose one:"

$survey.Range("B11").Value = "select_one yes_no"
$survey.Range("C11").Value = "inline"
$survey.Range("E11").Value = "i2"
$survey.Range("F11").Value = "Choose one:"

$survey.Range("B12").Value = "select_one yes_no"
$survey.Range("C12").Value = "inline"
$survey.Range("E12").Value = "i3"
$survey.Range("F12").Value = "Choose one:"

$survey.Range("B13").Value = "end screen"

# ---------------------------------------------------------------
# Sheet "choices" (sheet 2): add a "label" column, and add the
# yes/no choice list used by the new select_one yes_no question.
# ---------------------------------------------------------------
$choices = $wb.Worksheets.Item(2)

$choices.Range("D1").Value = "label"

$choices.Range("A12").Value = "yes_no"
$choices.Range("B12").Value = "yes"
$choices.Range("D12").Value = "Yes"

$choices.Range("A13").Value = "yes_no"
$choices.Range("B13").Value = "no"
$choices.Range("D13").Value = "No"

# ---------------------------------------------------------------
# Sheet "queries" (sheet 3): drop the unused param.format/param.q
# columns and switch the odk_values row to a content:// URI.
# ---------------------------------------------------------------
$queries = $wb.Worksheets.Item(3)

$queries.Columns.Item(4).Delete()
$queries.Columns.Item(4).Delete()

$queries.Range("B5").Value = '"content://com.opendatakit.tables.ContentProvider/database_id/table_id/row_id"'

# ---------------------------------------------------------------
# Sheet "settings" (sheet 4): widen the "value" column.
# ---------------------------------------------------------------
$settings = $wb.Worksheets.Item(4)
$settings.Columns.Item(2).ColumnWidth = 70.46
